$wb = $excel.ActiveWorkbook

$oldText = "February 03 2026 17.29.55 EST"
$newText = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")

$a2 = [string]$wsAbout.Range("A2").Value2
$wsAbout.Range("A2").Value2 = $a2.Replace($oldText, $newText)

$a6 = [string]$wsAbout.Range("A6").Value2
$wsAbout.Range("A6").Value2 = $a6.Replace($oldText, $newText)

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
foreach ($r in 2..7) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19
    $val = [string]$cell.Value2
    $cell.Value2 = $val.Replace($oldText, $newText)
}
